$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Relocate the existing unique strings (keeps shared-string table order/index stable) ---
# E1 (water_need) must move out first because G1 will land on top of E1.
$ws.Range("E1").Cut($ws.Range("B7"))
$ws.Range("G1").Cut($ws.Range("E1"))
$ws.Range("F1").Cut($ws.Range("C7"))
$ws.Range("H1").Cut($ws.Range("B13"))
$ws.Range("I1").Cut($ws.Range("C13"))
$ws.Range("J1").Cut($ws.Range("D13"))
$ws.Range("F2").Cut($ws.Range("C8"))
$ws.Range("F4").Cut($ws.Range("C10"))

# --- New header string (appended to the shared-string table) ---
$ws.Range("D7").Value = "water_received"

# --- Repeat the row-label strings in the two new mini tables ---
$ws.Range("A8").Value = "Empty"
$ws.Range("A9").Value = "Town"
$ws.Range("A10").Value = "River"
$ws.Range("A11").Value = "Tree"

$ws.Range("A14").Value = "Empty"
$ws.Range("A15").Value = "Town"
$ws.Range("A16").Value = "River"
$ws.Range("A17").Value = "Tree"

# --- Wipe whatever used to live outside the new A1:E17 footprint ---
$ws.Range("F1:J17").Clear()
$ws.Range("A18:J30").Clear()

# --- Apply the "Satisfaisant" (green) cell style to the highlighted cells that
#     don't already carry it from their old location. Do this BEFORE adding the
#     borders below, since (re)applying a named Style resets direct formatting. ---
$ws.Range("B11").Style = "Satisfaisant"
$ws.Range("D11").Style = "Satisfaisant"
$ws.Range("C14").Style = "Satisfaisant"
$ws.Range("C15").Style = "Satisfaisant"
$ws.Range("B17:D17").Style = "Satisfaisant"

# --- Thin borders around each of the three little tables ---
$ws.Range("A1:E5").Borders.LineStyle = 1
$ws.Range("A7:D11").Borders.LineStyle = 1
$ws.Range("A13:D17").Borders.LineStyle = 1

# --- Column A is narrower than the rest ---
$ws.Range("A1").EntireColumn.ColumnWidth = 9.140625

# --- Selection shown when the file was last saved ---
$ws.Range("F13").Select()
